{"js": "// The document contains a single paragraph reading \"Version 1.\" (with a\n// \"_GoBack\" bookmark sitting right after the \"1.\"). The edit bumps the\n// version number so the paragraph reads \"Version 2.\" instead.\nconst body = context.document.body;\n\n// Find the \"1.\" right after \"Version \" and swap it for \"2.\" in place so\n// the rest of the paragraph (the bookmark, any other runs/formatting)\n// is left untouched.\nconst results = body.search(\"1.\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"2.\", \"Replace\");\n} else {\n  // Fallback: if for some reason the expected \"1.\" substring is not\n  // found (e.g. the document text already changed), replace the whole\n  // paragraph text wholesale so the intent of the edit still lands.\n  const paragraphs = body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n  if (paragraphs.items.length > 0) {\n    paragraphs.items[0].load(\"text\");\n    await context.sync();\n    const newText = paragraphs.items[0].text.replace(\"1.\", \"2.\");\n    paragraphs.items[0].insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains a single paragraph reading \"Version 1.\" (with a\n# \"_GoBack\" bookmark sitting right after the \"1.\"). The edit bumps the\n# version number so the paragraph reads \"Version 2.\" instead.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"1.\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"2.\"\n$found = $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n\nif (-not $found) {\n    # Fallback: if the expected \"1.\" substring isn't present (e.g. the\n    # document text already changed), replace the whole paragraph text.\n    $para = $d.Paragraphs.Item(1).Range\n    $para.Text = $para.Text -replace \"1\\.\", \"2.\"\n}\n"}
